# Generate Report for handback
# Updates the zh-cn and de-de localization-status sheets to reflect that the
# files have been handed back (are now in sync with en-US): refreshes the
# status text, records the "Latest Target File" / "Latest Handback File"
# hyperlinks, and stamps the "Latest Handback DateTime" for each tracked
# source file.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet - mirrors the same status text for both locales
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value2 = $statusText
$wsOverview.Range("C2").Value2 = $statusText
$wsOverview.Range("B3").Value2 = $statusText
$wsOverview.Range("C3").Value2 = $statusText

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value2 = $statusText
$wsZh.Range("B3").Value2 = $statusText

# Row 2 -> 403c8a17-a961-4f06-8f77-0acf83e9e6b4 (.md / .zh-cn.xlf)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8acf5aca0c90e9d02d64f49418cbc7711a438868/e2e/403c8a17-a961-4f06-8f77-0acf83e9e6b4.md",
    "",
    "",
    "403c8a17-a961-4f06-8f77-0acf83e9e6b4.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b65ecebab9b5fab892cd9c49324513e048b08cc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/403c8a17-a961-4f06-8f77-0acf83e9e6b4.87edcfba4b6b414c0cc1968245d38807e519136a.zh-cn.xlf",
    "",
    "",
    "403c8a17-a961-4f06-8f77-0acf83e9e6b4.87edcfba4b6b414c0cc1968245d38807e519136a.zh-cn.xlf"
) | Out-Null
$wsZh.Range("G2").Value2 = "2016-01-18 06:37:45"

# Row 3 -> 606237ae-7f0d-486e-b353-7c2732e49989 (.md / .zh-cn.xlf)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8acf5aca0c90e9d02d64f49418cbc7711a438868/e2e/606237ae-7f0d-486e-b353-7c2732e49989.md",
    "",
    "",
    "606237ae-7f0d-486e-b353-7c2732e49989.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b65ecebab9b5fab892cd9c49324513e048b08cc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/606237ae-7f0d-486e-b353-7c2732e49989.7e32aa246cd8c7388eb5b50f4c54ff80db5fa22e.zh-cn.xlf",
    "",
    "",
    "606237ae-7f0d-486e-b353-7c2732e49989.7e32aa246cd8c7388eb5b50f4c54ff80db5fa22e.zh-cn.xlf"
) | Out-Null
$wsZh.Range("G3").Value2 = "2016-01-18 06:37:45"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value2 = $statusText
$wsDe.Range("B3").Value2 = $statusText

# Row 2 -> 403c8a17-a961-4f06-8f77-0acf83e9e6b4 (.md / .de-de.xlf)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8acf5aca0c90e9d02d64f49418cbc7711a438868/e2e/403c8a17-a961-4f06-8f77-0acf83e9e6b4.md",
    "",
    "",
    "403c8a17-a961-4f06-8f77-0acf83e9e6b4.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/243570438d61382568bd1827d15b69db3ffe8e47/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/403c8a17-a961-4f06-8f77-0acf83e9e6b4.87edcfba4b6b414c0cc1968245d38807e519136a.de-de.xlf",
    "",
    "",
    "403c8a17-a961-4f06-8f77-0acf83e9e6b4.87edcfba4b6b414c0cc1968245d38807e519136a.de-de.xlf"
) | Out-Null
$wsDe.Range("G2").Value2 = "2016-01-18 06:38:02"

# Row 3 -> 606237ae-7f0d-486e-b353-7c2732e49989 (.md / .de-de.xlf)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8acf5aca0c90e9d02d64f49418cbc7711a438868/e2e/606237ae-7f0d-486e-b353-7c2732e49989.md",
    "",
    "",
    "606237ae-7f0d-486e-b353-7c2732e49989.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/243570438d61382568bd1827d15b69db3ffe8e47/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/606237ae-7f0d-486e-b353-7c2732e49989.7e32aa246cd8c7388eb5b50f4c54ff80db5fa22e.de-de.xlf",
    "",
    "",
    "606237ae-7f0d-486e-b353-7c2732e49989.7e32aa246cd8c7388eb5b50f4c54ff80db5fa22e.de-de.xlf"
) | Out-Null
$wsDe.Range("G3").Value2 = "2016-01-18 06:38:02"
